# "Added zine, movie reviews, art section"
#
# Appends 16 new movie-review rows (48-63) to the bottom of Sheet1's review
# table. Sheet2 ("toSee") is left untouched - its cell text is unchanged by
# the commit (only the shared-string indices backing it shift, which is an
# OOXML storage detail, not a content change).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each hashtable is one new row; keys are column letters (A-J) mapped below.
# Omitted keys mean the source cell is blank, exactly like the diff.
$newRows = @(
    @{ row = 48; A = "Get Out"; B = 2017; C = "Jordan Peele"; D = "1h 44m "; E = "Horror, Comedy"; F = 8; G = 4; H = 6; J = "https://en.wikipedia.org/wiki/Get_Out" }
    @{ row = 49; A = "Raw"; B = 2016; C = "Julia Ducournau"; D = "1h 33m"; E = "Horror, Drama"; F = 8; G = 8; H = 9; I = 7; J = "https://en.wikipedia.org/wiki/Raw_(film)" }
    @{ row = 50; A = "Thoroughbreds"; B = 2017; C = "Cory Finley"; D = "1h 32m"; E = "Drama, Horror, Crime"; F = 7; G = 7; H = 9; J = "https://en.wikipedia.org/wiki/Thoroughbreds_(2017_film)" }
    @{ row = 51; A = "The Truman Show"; B = 1998; E = "Drama"; G = 5; H = 8 }
    @{ row = 52; A = "Persona"; B = 1966; E = "Drama, Existential Horror"; F = 2; G = 10; H = 8 }
    @{ row = 53; A = "The Wailing"; B = 2016; E = "Horror, Drama"; G = 9; H = 9 }
    @{ row = 54; A = "Peeping Tom"; B = 1960 }
    @{ row = 55; A = "Akira"; G = 9; H = 9 }
    @{ row = 56; A = "The Adventures of Baron Munchausen"; B = 1988 }
    @{ row = 57; A = "Melancholia" }
    @{ row = 58; A = "Mulholland Drive" }
    @{ row = 59; A = "Children of Men" }
    @{ row = 60; A = "Inside Out" }
    @{ row = 61; A = "Arrival" }
    @{ row = 62; A = "Mandy" }
    @{ row = 63; A = "Eastern Promises"; J = "https://en.wikipedia.org/wiki/Eastern_Promises" }
)

$colNum = @{ A = 1; B = 2; C = 3; D = 4; E = 5; F = 6; G = 7; H = 8; I = 9; J = 10 }

foreach ($r in $newRows) {
    $rowNum = $r.row
    foreach ($col in @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")) {
        if ($r.ContainsKey($col)) {
            $ws.Cells.Item($rowNum, $colNum[$col]).Value = $r[$col]
        }
    }
}

# Leave the selection where the author's editing session ended up (the
# dimension now runs through row 63, L47 -> L63).
$ws.Range("O52").Select()
